# The deck originally carries the "Integral" theme (ppt/theme/theme1.xml,
# applied via the slide master to every slide) and the generic default
# "Office Theme" palette tucked away as ppt/theme/theme2.xml (only used by
# the notes master). The edit swaps the two: the slide master's theme
# becomes the plain default Office palette and the former "Integral"
# palette moves to the notes-master theme.
#
# This COM-interop host does not expose a way to re-point a master at a
# different theme part, nor to rename a Theme/ColorScheme, nor to reach
# the notes-master's theme at all (Master.ColorScheme always resolves to
# the slide master's theme regardless of which master object it is read
# from). The one reliable, non-destructive lever available is
# Slide.ThemeColorScheme, which writes straight through to the slide
# master's theme part (ppt/theme/theme1.xml) without clobbering the
# <a:clrScheme name="..."/> attribute the way the legacy
# Slide.ColorScheme / Master.ColorScheme collections do. Use it to push
# every one of the 12 theme colors over to the standard Office palette.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Hex-ToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order matches msoThemeColor*: 1=dk1 2=lt1 3=dk2 4=lt2 5-10=accent1-6
# 11=hlink 12=folHlink. Values below are the stock "Office" theme palette.
$officeTheme = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = Hex-ToRGB($officeTheme[$i - 1])
}
